# Updated cryptos list on Mon Apr 10 08:41:47 UTC 2023 with GitHub Actions
# Refresh the live price / 1h-volume-change snapshot in the "cryptos" sheet.
# For numeric-looking price strings in column D, a leading apostrophe is used
# so Excel stores them as literal text (matching the source data, which can
# include values like "1.220" or "0.06730" whose trailing zeros must survive).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.533.63"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "1.874.67"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("E4").Value = "  -1.80%  "
$ws.Range("D5").Value = "'315.67"
$ws.Range("E5").Value = "  -0.56%  "
$ws.Range("E7").Value = "  -1.04%  "
$ws.Range("D8").Value = "'0.3906"
$ws.Range("E8").Value = "  -1.02%  "
$ws.Range("D9").Value = "'0.08362"
$ws.Range("E9").Value = "  +0.34%  "
$ws.Range("D10").Value = "'1.106"
$ws.Range("E10").Value = "  -1.05%  "
$ws.Range("D11").Value = "'41.82"
$ws.Range("E11").Value = "  -0.73%  "
$ws.Range("D12").Value = "'6.226"
$ws.Range("D13").Value = "1.876.80"
$ws.Range("E13").Value = "  +0.29%  "
$ws.Range("D14").Value = "'20.44"
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("D15").Value = "'7.286"
$ws.Range("E16").Value = "  -1.85%  "
$ws.Range("E17").Value = "  -0.46%  "
$ws.Range("D18").Value = "'91.26"
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("D19").Value = "'0.06730"
$ws.Range("E19").Value = "  -0.39%  "
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("E21").Value = "  -1.53%  "
$ws.Range("D22").Value = "'5.917"
$ws.Range("E22").Value = "  -0.82%  "
$ws.Range("D23").Value = "28.569.17"
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").Value = "'11.17"
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("D25").Value = "'2.225"
$ws.Range("E25").Value = "  -1.90%  "
$ws.Range("D26").Value = "2.094.62"
$ws.Range("E26").Value = "  +0.49%  "
$ws.Range("D27").Value = "'161.55"
$ws.Range("E27").Value = "  -0.25%  "
$ws.Range("E28").Value = "  -0.74%  "
$ws.Range("D29").Value = "'2.419"
$ws.Range("E29").Value = "  +2.13%  "
$ws.Range("D30").Value = "'126.50"
$ws.Range("E30").Value = "  -0.78%  "
$ws.Range("D31").Value = "'0.1042"
$ws.Range("E31").Value = "  -1.09%  "
$ws.Range("D32").Value = "'1.042"
$ws.Range("E32").Value = "  +0.68%  "
$ws.Range("D33").Value = "'5.744"
$ws.Range("E33").Value = "  -1.88%  "
$ws.Range("D34").Value = "'3.615"
$ws.Range("E34").Value = "  -1.27%  "
$ws.Range("D35").Value = "'0.02453"
$ws.Range("E35").Value = "  +0.55%  "
$ws.Range("D36").Value = "'0.06556"
$ws.Range("E36").Value = "  +1.00%  "
$ws.Range("D37").Value = "'8.932"
$ws.Range("E37").Value = "  -2.36%  "
$ws.Range("E38").Value = "  -0.63%  "
$ws.Range("D39").Value = "'5.028"
$ws.Range("E39").Value = "  +0.92%  "
$ws.Range("E40").Value = "  -0.41%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.6387"
$ws.Range("E41").Value = "  -1.09%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'1.234"
$ws.Range("E42").Value = "  -1.36%  "
$ws.Range("D43").Value = "'11.11"
$ws.Range("E43").Value = "  -0.65%  "
$ws.Range("D45").Value = "'0.6010"
$ws.Range("E45").Value = "  -0.37%  "
$ws.Range("D46").Value = "'13.02"
$ws.Range("E46").Value = "  +0.80%  "
$ws.Range("D47").Value = "'3.691"
$ws.Range("E47").Value = "  -0.79%  "
$ws.Range("D48").Value = "'2.004"
$ws.Range("E48").Value = "  +0.55%  "
$ws.Range("D49").Value = "'1.220"
$ws.Range("E49").Value = "  +0.39%  "
$ws.Range("D50").Value = "'122.09"
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("D51").Value = "'1.143"
$ws.Range("E51").Value = "  -9.80%  "

